$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current state: C2 = "192.168.1.113" (text-formatted), E2 = "192.168.0.24" (general format)
# Target state:  C2 = "127.0.0.1"    (text-formatted), E2 = "192.168.0.24" (text-formatted)

# Give E2 the same (text) number format C2 already has; E2's own value stays the same.
$ws.Range("E2").NumberFormat = $ws.Range("C2").NumberFormat

# C2 gets the brand-new IP address, replacing "192.168.1.113" (which then drops
# out of the workbook entirely, since no cell references it any more).
$ws.Range("C2").Value2 = "127.0.0.1"

$ws.Range("E2").Select()
